# "fix urban & farmland area"
#
# 1. Economy sheet: insert a new row (row 5) "Urban area per person
#    (ha=100m*100m)" computed from the existing "Urban area per person (m²)"
#    row (row 4) by converting m² -> ha (i.e. /100/100). This pushes the old
#    rows 5-7 (Urban construction height / Farmer to eater ratio / Fisher to
#    eater ratio) down to rows 6-8, and updates "Urban construction height"
#    with corrected values.
# 2. Land Use sheet: the Min/Max-area-per-citizen formulas (D3/D4) are
#    updated to use the new "Urban area per person (ha)" row together with
#    "Urban construction height" instead of the old (incorrect) formula that
#    divided square metres by 100 instead of converting properly to hectares.
# 3. Citizen Burg Modifiers sheet: Farmer now also gets a Naval modifier
#    (-10).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Economy sheet
# ---------------------------------------------------------------------
$economy = $wb.Worksheets.Item("Economy")

# Insert a new row above the old "Urban construction height" row (row 5)
$economy.Rows("5:5").Insert()

$economy.Range("A5").Value = "Urban area per person (ha=100m*100m)"

$economy.Range("B5").Formula = "=B4/100/100"
$economy.Range("C5").Formula = "=C4/100/100"
$economy.Range("D5").Formula = "=D4/100/100"
$economy.Range("E5").Formula = "=E4/100/100"

# Corrected "Urban construction height" values (now row 6)
$economy.Range("B6").Value = 2.0
$economy.Range("C6").Value = 4.0
$economy.Range("D6").Value = 3.0
$economy.Range("E6").Value = 10.0

# ---------------------------------------------------------------------
# 2. Land Use sheet
# ---------------------------------------------------------------------
$landUse = $wb.Worksheets.Item("Land Use")
$landUse.Range("D3").Formula = "=B3*Economy!B5/Economy!C6"
$landUse.Range("D4").Formula = "=C4*Economy!C5/Economy!B6"

# ---------------------------------------------------------------------
# 3. Citizen Burg Modifiers sheet
# ---------------------------------------------------------------------
$burgMod = $wb.Worksheets.Item("Citizen Burg Modifiers")
$burgMod.Range("G2").Copy()
$burgMod.Range("D2").PasteSpecial(-4122)
$burgMod.Range("D2").Value = -10.0

$excel.CalculateFullRebuild()
